$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 462.16666
$ws.Range("J9").Value = 574.3333
$ws.Range("L9").Value = 574.3333
$ws.Range("N9").Value = -912.3333
$ws.Range("H18").Value = 2245.2727
$ws.Range("I18").Value = 1969.8
$ws.Range("J18").Value = 5000
$ws.Range("K18").Value = 1969.8
$ws.Range("L18").Value = 5000
$ws.Range("M18").Value = -1685.8
$ws.Range("N18").Value = -5568
$ws.Range("H33").Value = 164.84616
$ws.Range("I33").Value = 164.84616
$ws.Range("K33").Value = 164.84616
$ws.Range("M33").Value = 64.15384
$ws.Range("H111").Value = 3598.4
$ws.Range("I111").Value = 3598.4
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 10795.2
$ws.Range("L111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -7728.200000000001
$ws.Range("H132").Value = 1805.5333
$ws.Range("I132").Value = 1937.5385
$ws.Range("K132").Value = 5812.6155
$ws.Range("M132").Value = -3282.6155
$ws.Range("H137").Value = 2435.0435
$ws.Range("I137").Value = 1156
$ws.Range("K137").Value = 3468
$ws.Range("M137").Value = -918

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4059.4443
$ws.Range("J61").Value = 4104.8
$ws.Range("L61").Value = 4104.8
$ws.Range("N61").Value = -4528.8
$ws.Range("H104").Value = 5000
$ws.Range("J104").Value = 5000
$ws.Range("L104").Value = 5000
$ws.Range("N104").Value = -11988
$ws.Range("H122").Value = 402602.9
$ws.Range("I122").Value = 715755.2
$ws.Range("K122").Value = 2147265.6
$ws.Range("M122").Value = -2144815.6
$ws.Range("H136").Value = 4059.4443
$ws.Range("J136").Value = 4104.8
$ws.Range("L136").Value = 12314.4
$ws.Range("N136").Value = -17414.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 997.3333
$ws.Range("I94").Value = 1277.4667
$ws.Range("J94").Value = 530.44446
$ws.Range("K94").Value = 1277.4667
$ws.Range("L94").Value = 530.44446
$ws.Range("M94").Value = -826.4666999999999
$ws.Range("N94").Value = -1432.44446
$ws.Range("H134").Value = 1935.742
$ws.Range("I134").Value = 1642.6923
$ws.Range("K134").Value = 4928.0769
$ws.Range("M134").Value = -2393.0769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 156.86957
$ws.Range("I7").Value = 124.2381
$ws.Range("K7").Value = 124.2381
$ws.Range("M7").Value = -11.2381
$ws.Range("H68").Value = 38333
$ws.Range("J68").Value = 38333
$ws.Range("L68").Value = 38333
$ws.Range("N68").Value = -39831
$ws.Range("H71").Value = 38333
$ws.Range("J71").Value = 38333
$ws.Range("L71").Value = 114999
$ws.Range("N71").Value = -122487
$ws.Range("H82").Value = 3000
$ws.Range("J82").Value = 3000
$ws.Range("L82").Value = 3000
$ws.Range("N82").Value = -3722
$ws.Range("H85").Value = 3000
$ws.Range("J85").Value = 3000
$ws.Range("L85").Value = 3000
$ws.Range("N85").Value = -5496
$ws.Range("H86").Value = 8667.333000000001
$ws.Range("I86").Value = 6501
$ws.Range("K86").Value = 6501
$ws.Range("M86").Value = -5378
$ws.Range("H89").Value = 8667.333000000001
$ws.Range("I89").Value = 6501
$ws.Range("K89").Value = 32505
$ws.Range("M89").Value = -26889
$ws.Range("H122").Value = 2087.6924
$ws.Range("I122").Value = 2134.5833
$ws.Range("K122").Value = 6403.749899999999
$ws.Range("M122").Value = -3953.749899999999
$ws.Range("H132").Value = 3612
$ws.Range("I132").Value = 3612
$ws.Range("K132").Value = 10836
$ws.Range("M132").Value = -8306
$ws.Range("H134").Value = 3267.5625
$ws.Range("I134").Value = 2481.9092
$ws.Range("K134").Value = 7445.7276
$ws.Range("M134").Value = -4910.7276
$ws.Range("H141").Value = 192183
$ws.Range("I141").Value = 113275
$ws.Range("K141").Value = 113275
$ws.Range("M141").Value = -108095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 59658028
$ws.Range("J4").Value = 800
$ws.Range("L4").Value = 2400
$ws.Range("N4").Value = -2624
$ws.Range("H12").Value = 853.5833
$ws.Range("I12").Value = 833
$ws.Range("K12").Value = 2499
$ws.Range("M12").Value = -2326
$ws.Range("H14").Value = 614.1111
$ws.Range("I14").Value = 614.1111
$ws.Range("K14").Value = 1842.3333
$ws.Range("M14").Value = -1669.3333
$ws.Range("H33").Value = 1556.5
$ws.Range("J33").Value = 938
$ws.Range("L33").Value = 5628
$ws.Range("N33").Value = -6194
$ws.Range("H34").Value = 2244
$ws.Range("J34").Value = 4000
$ws.Range("L34").Value = 12000
$ws.Range("N34").Value = -12168
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H42").Value = 725
$ws.Range("I42").Value = 700
$ws.Range("J42").Value = 800
$ws.Range("K42").Value = 2100
$ws.Range("L42").Value = 2400
$ws.Range("M42").Value = -1566
$ws.Range("N42").Value = -3468
$ws.Range("H117").Value = 2479.2222
$ws.Range("J117").Value = 3343.8333
$ws.Range("L117").Value = 10031.4999
$ws.Range("N117").Value = -16915.4999
$ws.Range("H132").Value = 1899.75
$ws.Range("J132").Value = 2299.5
$ws.Range("L132").Value = 20695.5
$ws.Range("N132").Value = -25755.5
$ws.Range("H136").Value = 19125
$ws.Range("I136").Value = 1500
$ws.Range("K136").Value = 4500
$ws.Range("M136").Value = 600

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 249.55
$ws.Range("I2").Value = 62.2
$ws.Range("J2").Value = 436.9
$ws.Range("K2").Value = 62.2
$ws.Range("L2").Value = 436.9
$ws.Range("M2").Value = 50.8
$ws.Range("N2").Value = -662.9
$ws.Range("H107").Value = 626.9091
$ws.Range("J107").Value = 782.8570999999999
$ws.Range("L107").Value = 782.8570999999999
$ws.Range("N107").Value = -4622.8571
$ws.Range("H122").Value = 61647.53
$ws.Range("I122").Value = 2600.6428
$ws.Range("K122").Value = 7801.928400000001
$ws.Range("M122").Value = -5351.928400000001
$ws.Range("H132").Value = 3679.625
$ws.Range("I132").Value = 2300
$ws.Range("J132").Value = 7818.5
$ws.Range("K132").Value = 6900
$ws.Range("L132").Value = 23455.5
$ws.Range("M132").Value = -4370
$ws.Range("N132").Value = -28515.5
$ws.Range("H138").Value = 31578.908
$ws.Range("J138").Value = 50429
$ws.Range("L138").Value = 50429
$ws.Range("N138").Value = -60709

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 13999.5
$ws.Range("I22").Value = 2499
$ws.Range("J22").Value = 17833
$ws.Range("K22").Value = 2499
$ws.Range("L22").Value = 17833
$ws.Range("M22").Value = -2204
$ws.Range("N22").Value = -18423
$ws.Range("H27").Value = 13999.5
$ws.Range("I27").Value = 2499
$ws.Range("J27").Value = 17833
$ws.Range("K27").Value = 2499
$ws.Range("L27").Value = 17833
$ws.Range("M27").Value = -2392
$ws.Range("N27").Value = -18047
$ws.Range("H68").Value = 1958
$ws.Range("J68").Value = 1790
$ws.Range("L68").Value = 1790
$ws.Range("N68").Value = -3288
$ws.Range("H71").Value = 1958
$ws.Range("J71").Value = 1790
$ws.Range("L71").Value = 8950
$ws.Range("N71").Value = -16438
$ws.Range("H82").Value = 2587.6206
$ws.Range("J82").Value = 1559.2
$ws.Range("L82").Value = 1559.2
$ws.Range("N82").Value = -2281.2
$ws.Range("H85").Value = 2587.6206
$ws.Range("J85").Value = 1559.2
$ws.Range("L85").Value = 1559.2
$ws.Range("N85").Value = -4055.2
$ws.Range("H106").Value = 20126.5
$ws.Range("J106").Value = 20126.5
$ws.Range("L106").Value = 20126.5
$ws.Range("N106").Value = -22650.5
$ws.Range("H136").Value = 3785
$ws.Range("I136").Value = 3399.2
$ws.Range("K136").Value = 10197.6
$ws.Range("M136").Value = -7647.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1575.4
$ws.Range("I136").Value = 1195.25
$ws.Range("J136").Value = 4996.75
$ws.Range("K136").Value = 3585.75
$ws.Range("L136").Value = 14990.25
$ws.Range("M136").Value = -1035.75
$ws.Range("N136").Value = -20090.25
